# Update workbook data as of 2025-11-05 12:33 run.
$wb = $excel.ActiveWorkbook

# --- 1. Update "Last Updated" timestamp on the Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "05 Nov 2025, 12:32 PM"

# --- 2. Update the "Stock List" sheet: a new stock (CAPTRU-RE1) now leads the
#        list, pushing every existing row down by one and dropping the last
#        row (TRAVELFOOD) off the bottom of the table. ---
$ws = $wb.Worksheets.Item("Stock List")

# Insert a fresh row at the top of the data (row 2, just below the header)
$ws.Rows.Item(2).Insert()

# The freshly inserted row inherits formatting from the row above (the bold
# header); reset it so it matches the rest of the plain data rows.
$ws.Range("A2:H2").ClearFormats()

# Populate the new row with the new stock's data. Icon / 52W High /
# Distance 52W High / Market Cap keep the same "template" values the table
# already used for this position, so copy them from the row just below.
$ws.Range("A2").Value = $ws.Range("A3").Value2
$ws.Range("B2").Value = "CAPTRU-RE1"
$ws.Range("C2").Value = "CAPTRU-RE1"
$ws.Range("D2").Value = 5.67
$ws.Range("E2").Value = -11.9565
$ws.Range("F2").Value = $ws.Range("F3").Value2
$ws.Range("G2").Value = $ws.Range("G3").Value2
$ws.Range("H2").Value = $ws.Range("H3").Value2

# Remove the row that got pushed past the end of the original table
# (previously row 76 / TRAVELFOOD, now shifted to row 77).
$ws.Rows.Item(77).Delete()
